$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 85, shifting existing rows 85-187 down to 86-188
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with data
$ws.Cells.Item(85, 1).Value = 6
$ws.Cells.Item(85, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(85, 3).Value = "Metropolitana"
$ws.Cells.Item(85, 4).Value = 44546
$ws.Cells.Item(85, 5).Value = 13
$ws.Cells.Item(85, 6).Value = 100112022
$ws.Cells.Item(85, 7).Value = "Arveja Verde"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 700
$ws.Cells.Item(85, 11).Value = 9000
$ws.Cells.Item(85, 12).Value = 10000
$ws.Cells.Item(85, 13).Value = 9429
$ws.Cells.Item(85, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(85, 15).Value = "Carahue"
$ws.Cells.Item(85, 16).Value = 377
$ws.Cells.Item(85, 17).Value = 25
$ws.Cells.Item(85, 18).Value = "Hortaliza"
